$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (matching source data) for Price-column cells whose new
# value would otherwise be auto-recognized as a plain number by Excel.
$ws.Cells.Item(2, 4).Value = '68.319.20'
$ws.Cells.Item(2, 5).Value = '  +0.93%  '
$ws.Cells.Item(3, 4).Value = '3.351.57'
$ws.Cells.Item(3, 5).Value = '  +0.59%  '
$ws.Cells.Item(4, 5).Value = '  -0.03%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '583.55'
$ws.Cells.Item(5, 5).Value = '  +0.28%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '177.39'
$ws.Cells.Item(6, 5).Value = '  +1.05%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.999'
$ws.Cells.Item(7, 5).Value = '  +0.07%  '
$ws.Cells.Item(8, 5).Value = '  +0.41%  '
$ws.Cells.Item(9, 5).Value = '  +2.83%  '
$ws.Cells.Item(10, 5).Value = '  +0.92%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '48.10'
$ws.Cells.Item(11, 5).Value = '  +5.46%  '
$ws.Cells.Item(12, 5).Value = '  +1.25%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '692.06'
$ws.Cells.Item(13, 5).Value = '  +3.90%  '
$ws.Cells.Item(14, 4).Value = '3.886.10'
$ws.Cells.Item(14, 5).Value = '  +0.36%  '
$ws.Cells.Item(15, 5).Value = '  +0.16%  '
$ws.Cells.Item(16, 4).Value = '68.349.68'
$ws.Cells.Item(16, 5).Value = '  +0.65%  '
$ws.Cells.Item(17, 5).Value = '  +1.25%  '
$ws.Cells.Item(18, 4).Value = '3.360.08'
$ws.Cells.Item(18, 5).Value = '  +0.92%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '17.47'
$ws.Cells.Item(19, 5).Value = '  -0.04%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '11.20'
$ws.Cells.Item(20, 5).Value = '  +2.29%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '0.896'
$ws.Cells.Item(21, 5).Value = '  +0.55%  '
$ws.Cells.Item(22, 5).Value = '  +0.81%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '17.01'
$ws.Cells.Item(23, 5).Value = '  -0.54%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '100.04'
$ws.Cells.Item(24, 5).Value = '  +0.86%  '
$ws.Cells.Item(25, 5).Value = '  +1.54%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '2.70'
$ws.Cells.Item(26, 5).Value = '  +0.61%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '9.55'
$ws.Cells.Item(27, 5).Value = '  +2.77%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '33.04'
$ws.Cells.Item(28, 5).Value = '  -1.99%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '8.50'
$ws.Cells.Item(29, 5).Value = '  +0.59%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '6.94'
$ws.Cells.Item(30, 5).Value = '  -6.94%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '559.11'
$ws.Cells.Item(31, 5).Value = '  -5.54%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '11.07'
$ws.Cells.Item(32, 5).Value = '  +1.00%  '
$ws.Cells.Item(33, 5).Value = '  +1.09%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '57.60'
$ws.Cells.Item(34, 5).Value = '  +1.37%  '
$ws.Cells.Item(35, 5).Value = '  +0.16%  '
$ws.Cells.Item(36, 4).Value = '3.704.18'
$ws.Cells.Item(37, 5).Value = '  +0.37%  '
$ws.Cells.Item(38, 5).Value = '  +3.75%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '34.67'
$ws.Cells.Item(39, 5).Value = '  +3.43%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '3.18'
$ws.Cells.Item(40, 5).Value = '  +2.08%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '2.62'
$ws.Cells.Item(41, 5).Value = '  -0.54%  '
$ws.Cells.Item(42, 4).Value = '0.0₃0673'
$ws.Cells.Item(42, 5).Value = '  +1.24%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.336'
$ws.Cells.Item(43, 5).Value = '  +0.75%  '
$ws.Cells.Item(44, 5).Value = '  +1.30%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.0412'
$ws.Cells.Item(45, 5).Value = '  +1.17%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '2.66'
$ws.Cells.Item(46, 5).Value = '  +2.34%  '
$ws.Cells.Item(47, 5).Value = '  +0.61%  '
$ws.Cells.Item(48, 5).Value = '  -0.20%  '
$ws.Cells.Item(49, 5).Value = '  -0.34%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '130.86'
$ws.Cells.Item(50, 5).Value = '  +2.96%  '
$ws.Cells.Item(51, 5).Value = '  -0.29%  '
